# TCTC_02 with extent reports
# Adds a new "Sheet2" (Address details) after the existing "Sheet1",
# populates it with Country/City/Address/Postal Code/Phone/Fax data,
# and makes it the active sheet/tab.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# Insert the new worksheet right after Sheet1.
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Sheet2"

# Enter the header / data cells in the same left-to-right, top-to-bottom
# order they were authored in (column B.. first, then column A), so the
# shared-strings table ends up populated in the original order.
$ws2.Range("B1").Value = "City"
$ws2.Range("C1").Value = "Address 1"
$ws2.Range("D1").Value = "Address 2"
$ws2.Range("E1").Value = "Postal Code"
$ws2.Range("F1").Value = "Phone number"
$ws2.Range("G1").Value = "FaxNumber"

$ws2.Range("B2").Value = "Delhi"
$ws2.Range("C2").Value = "Random #103 Delhi India"
$ws2.Range("D2").Value = "Random #103 Delhi India"

$ws2.Range("A1").Value = "Country"
$ws2.Range("A2").Value = "India"

$ws2.Range("E2").Value = 123456
$ws2.Range("F2").Value = 9876543210
$ws2.Range("G2").Value = 55555555

# Column widths for the new sheet.
$ws2.Columns.Item(2).ColumnWidth = 14.44140625
$ws2.Columns.Item(3).ColumnWidth = 26.6640625
$ws2.Columns.Item(4).ColumnWidth = 17
$ws2.Columns.Item(5).ColumnWidth = 17.88671875
$ws2.Columns.Item(6).ColumnWidth = 17.5546875
$ws2.Columns.Item(7).ColumnWidth = 29.109375

# Page setup to match a single-page portrait printout.
$ws2.PageSetup.PaperSize = 9
$ws2.PageSetup.Orientation = 1

# Make Sheet2 the active/selected sheet, with A2 as the active cell and
# the view zoomed to 85%.
$ws2.Activate()
$ws2.Range("A2").Select()
$excel.ActiveWindow.Zoom = 85
